$d = $word.ActiveDocument

# Locate the "Presentación del proyecto" Heading2 paragraph. Right after it
# sits a duplicate screenshot (the same image already shown under the
# "Introducción" heading) that needs to be removed, per the commit:
# "eliminar capturas duplicadas ... (intro=featured, presentacion repetidas ...)".
$headingText = "Presentación del proyecto"

$count = $d.Paragraphs.Count
$targetIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Trim() -eq $headingText) {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -gt 0 -and $targetIndex -lt $d.Paragraphs.Count) {
    $nextPara = $d.Paragraphs.Item($targetIndex + 1)
    # Only remove it if it is indeed the picture paragraph (centered image,
    # no real text) so we don't accidentally delete something else.
    if ($nextPara.Range.InlineShapes.Count -gt 0 -and $nextPara.Range.Text.Trim() -eq "") {
        $nextPara.Range.Delete()
    }
}
